$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'26.843.98"
$c.Style = "Normal"
$c = $ws.Range("E2")
$c.Value = "'  -1.06%  "
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'1.873.21"
$c.Style = "Normal"
$c = $ws.Range("E3")
$c.Value = "'  -1.41%  "
$c.Style = "Normal"
$c = $ws.Range("E4")
$c.Value = "'  -0.30%  "
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'301.17"
$c.Style = "Normal"
$c = $ws.Range("E5")
$c.Value = "'  -1.93%  "
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'1.000"
$c.Style = "Normal"
$c = $ws.Range("E6")
$c.Value = "'  -0.26%  "
$c.Style = "Normal"
$c = $ws.Range("D7")
$c.Value = "'0.5342"
$c.Style = "Normal"
$c = $ws.Range("E7")
$c.Value = "'  +2.10%  "
$c.Style = "Normal"
$c = $ws.Range("D8")
$c.Value = "'0.3757"
$c.Style = "Normal"
$c = $ws.Range("E8")
$c.Value = "'  -1.28%  "
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.07179"
$c.Style = "Normal"
$c = $ws.Range("E9")
$c.Value = "'  -1.48%  "
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'21.64"
$c.Style = "Normal"
$c = $ws.Range("E10")
$c.Value = "'  +1.28%  "
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.8867"
$c.Style = "Normal"
$c = $ws.Range("E11")
$c.Value = "'  -1.87%  "
$c.Style = "Normal"
$c = $ws.Range("D12")
$c.Value = "'0.08112"
$c.Style = "Normal"
$c = $ws.Range("E12")
$c.Value = "'  -1.16%  "
$c.Style = "Normal"
$c = $ws.Range("D13")
$c.Value = "'1.865.94"
$c.Style = "Normal"
$c = $ws.Range("E13")
$c.Value = "'  -1.31%  "
$c.Style = "Normal"
$c = $ws.Range("D14")
$c.Value = "'93.08"
$c.Style = "Normal"
$c = $ws.Range("E14")
$c.Value = "'  -2.53%  "
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'5.275"
$c.Style = "Normal"
$c = $ws.Range("E15")
$c.Value = "'  -1.49%  "
$c.Style = "Normal"
$c = $ws.Range("E16")
$c.Value = "'  -0.25%  "
$c.Style = "Normal"
$c = $ws.Range("E17")
$c.Value = "'  +0.31%  "
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'0.000008541"
$c.Style = "Normal"
$c = $ws.Range("E18")
$c.Value = "'  -1.42%  "
$c.Style = "Normal"
$c = $ws.Range("E19")
$c.Value = "'  -0.21%  "
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'26.882.18"
$c.Style = "Normal"
$c = $ws.Range("E20")
$c.Value = "'  -1.08%  "
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'4.971"
$c.Style = "Normal"
$c = $ws.Range("E21")
$c.Value = "'  -3.12%  "
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'10.69"
$c.Style = "Normal"
$c = $ws.Range("E22")
$c.Value = "'  -0.95%  "
$c.Style = "Normal"
$c = $ws.Range("E23")
$c.Value = "'  -0.94%  "
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'147.16"
$c.Style = "Normal"
$c = $ws.Range("E24")
$c.Value = "'  -1.78%  "
$c.Style = "Normal"
$c = $ws.Range("E25")
$c.Value = "'  -3.10%  "
$c.Style = "Normal"
$c = $ws.Range("D26")
$c.Value = "'1.732"
$c.Style = "Normal"
$c = $ws.Range("E26")
$c.Value = "'  -0.45%  "
$c.Style = "Normal"
$c = $ws.Range("E27")
$c.Value = "'  -1.41%  "
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'114.37"
$c.Style = "Normal"
$c = $ws.Range("E28")
$c.Value = "'  -1.12%  "
$c.Style = "Normal"
$c = $ws.Range("D29")
$c.Value = "'4.747"
$c.Style = "Normal"
$c = $ws.Range("E29")
$c.Value = "'  -1.53%  "
$c.Style = "Normal"
$c = $ws.Range("E30")
$c.Value = "'  -6.61%  "
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'0.09134"
$c.Style = "Normal"
$c = $ws.Range("E31")
$c.Value = "'  -0.97%  "
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'0.7981"
$c.Style = "Normal"
$c = $ws.Range("E32")
$c.Value = "'  +0.72%  "
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'0.04989"
$c.Style = "Normal"
$c = $ws.Range("E33")
$c.Value = "'  -1.03%  "
$c.Style = "Normal"
$c = $ws.Range("D34")
$c.Value = "'2.988"
$c.Style = "Normal"
$c = $ws.Range("E34")
$c.Value = "'  +0.93%  "
$c.Style = "Normal"
$c = $ws.Range("D35")
$c.Value = "'1.171"
$c.Style = "Normal"
$c = $ws.Range("E35")
$c.Value = "'  -4.10%  "
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'0.5923"
$c.Style = "Normal"
$c = $ws.Range("E36")
$c.Value = "'  +3.39%  "
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'2.615"
$c.Style = "Normal"
$c = $ws.Range("E37")
$c.Value = "'  -0.57%  "
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'3.144"
$c.Style = "Normal"
$c = $ws.Range("E38")
$c.Value = "'  -6.44%  "
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.01949"
$c.Style = "Normal"
$c = $ws.Range("E39")
$c.Value = "'  -2.06%  "
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'1.068"
$c.Style = "Normal"
$c = $ws.Range("E40")
$c.Value = "'  -1.18%  "
$c.Style = "Normal"
$c = $ws.Range("D41")
$c.Value = "'6.661"
$c.Style = "Normal"
$c = $ws.Range("E41")
$c.Value = "'  +0.70%  "
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'8.910"
$c.Style = "Normal"
$c = $ws.Range("E42")
$c.Value = "'  -2.02%  "
$c.Style = "Normal"
$c = $ws.Range("D43")
$c.Value = "'115.70"
$c.Style = "Normal"
$c = $ws.Range("E43")
$c.Value = "'  -0.62%  "
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'0.5044"
$c.Style = "Normal"
$c = $ws.Range("E44")
$c.Value = "'  +3.01%  "
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'0.1494"
$c.Style = "Normal"
$c = $ws.Range("E45")
$c.Value = "'  -1.56%  "
$c.Style = "Normal"
$c = $ws.Range("E46")
$c.Value = "'  -0.36%  "
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'9.940"
$c.Style = "Normal"
$c = $ws.Range("E47")
$c.Value = "'  -1.88%  "
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'1.621"
$c.Style = "Normal"
$c = $ws.Range("E48")
$c.Value = "'  -0.86%  "
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'37.66"
$c.Style = "Normal"
$c = $ws.Range("E49")
$c.Value = "'  -2.19%  "
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'0.06028"
$c.Style = "Normal"
$c = $ws.Range("E50")
$c.Value = "'  +1.23%  "
$c.Style = "Normal"
$c = $ws.Range("D51")
$c.Value = "'62.18"
$c.Style = "Normal"
$c = $ws.Range("E51")
$c.Value = "'  -2.81%  "
$c.Style = "Normal"
